$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update CreatedAt timestamp in A1
$ws.Range("A1").Value = "CreatedAt: 2025-11-25T19:06:43"

# Update numeric data cells (columns W-Z, hours 21-24)
$ws.Range("X4").Value = 97.70999999999999
$ws.Range("Y4").Value = 697.36
$ws.Range("Z4").Value = 106.82
$ws.Range("W5").Value = -54.32
$ws.Range("X5").Value = -56.17
$ws.Range("Y5").Value = 0
$ws.Range("W6").Value = -15.11
$ws.Range("X6").Value = -11.61
$ws.Range("Y6").Value = -11.06
$ws.Range("Z6").Value = -4.7
$ws.Range("W8").Value = -24.43
$ws.Range("X8").Value = -14.36
$ws.Range("Y8").Value = 499.77
$ws.Range("W9").Value = 108.07
$ws.Range("X9").Value = 91.65000000000001
$ws.Range("Y9").Value = 698.11
$ws.Range("Z9").Value = 110.53
$ws.Range("W10").Value = -54.32
$ws.Range("X10").Value = -56.17
$ws.Range("Y10").Value = 0
$ws.Range("W11").Value = -22.04
$ws.Range("X11").Value = -17.68
$ws.Range("Y11").Value = -10.31
$ws.Range("Z11").Value = -0.99
$ws.Range("W13").Value = -24.43
$ws.Range("X13").Value = -14.36
$ws.Range("Y13").Value = 499.77
$ws.Range("W14").Value = 123.8
$ws.Range("X14").Value = 112
$ws.Range("Y14").Value = 698.11
$ws.Range("Z14").Value = 110.64
$ws.Range("W15").Value = -38.59
$ws.Range("X15").Value = -35.82
$ws.Range("W16").Value = -22.04
$ws.Range("X16").Value = -17.68
$ws.Range("Y16").Value = -10.31
$ws.Range("Z16").Value = -0.89
$ws.Range("W18").Value = -24.43
$ws.Range("X18").Value = -14.36
$ws.Range("Y18").Value = 499.77
$ws.Range("W19").Value = 115.3
$ws.Range("X19").Value = 96
$ws.Range("Y19").Value = 695.87
$ws.Range("Z19").Value = 106.52
$ws.Range("W20").Value = -51.89
$ws.Range("X20").Value = -56.17
$ws.Range("Y20").Value = 0
$ws.Range("W21").Value = -17.25
$ws.Range("X21").Value = -13.32
$ws.Range("Y21").Value = -12.55
$ws.Range("Z21").Value = -5.01
$ws.Range("W23").Value = -24.43
$ws.Range("X23").Value = -14.36
$ws.Range("Y23").Value = 499.77
$ws.Range("W24").Value = 112.87
$ws.Range("X24").Value = 96
$ws.Range("Y24").Value = 695.87
$ws.Range("Z24").Value = 106.52
$ws.Range("W25").Value = -54.32
$ws.Range("X25").Value = -56.17
$ws.Range("Y25").Value = 0
$ws.Range("W26").Value = -17.25
$ws.Range("X26").Value = -13.32
$ws.Range("Y26").Value = -12.55
$ws.Range("Z26").Value = -5.01
$ws.Range("W28").Value = -24.43
$ws.Range("X28").Value = -14.36
$ws.Range("Y28").Value = 499.77
$ws.Range("W29").Value = 113.39
$ws.Range("X29").Value = 94.62
$ws.Range("Y29").Value = 694.22
$ws.Range("Z29").Value = 106.01
$ws.Range("W30").Value = -51.89
$ws.Range("X30").Value = -56.17
$ws.Range("Y30").Value = 0
$ws.Range("W31").Value = -19.16
$ws.Range("X31").Value = -14.7
$ws.Range("Y31").Value = -14.2
$ws.Range("Z31").Value = -5.51
$ws.Range("W33").Value = -24.43
$ws.Range("X33").Value = -14.36
$ws.Range("Y33").Value = 499.77
$ws.Range("W34").Value = 121.82
$ws.Range("X34").Value = 105.81
$ws.Range("Y34").Value = 699.8200000000001
$ws.Range("Z34").Value = 113.11
$ws.Range("W35").Value = -38.59
$ws.Range("X35").Value = -35.82
$ws.Range("W36").Value = -24.03
$ws.Range("X36").Value = -19.27
$ws.Range("Y36").Value = -8.6
$ws.Range("Z36").Value = 1.58
$ws.Range("X37").Value = -4.6
$ws.Range("W38").Value = -24.43
$ws.Range("X38").Value = -14.36
$ws.Range("Y38").Value = 499.77
$ws.Range("X39").Value = 97.70999999999999
$ws.Range("Y39").Value = 697.36
$ws.Range("Z39").Value = 106.82
$ws.Range("W40").Value = -54.32
$ws.Range("X40").Value = -56.17
$ws.Range("Y40").Value = 0
$ws.Range("W41").Value = -15.11
$ws.Range("X41").Value = -11.61
$ws.Range("Y41").Value = -11.06
$ws.Range("Z41").Value = -4.7
$ws.Range("W43").Value = -24.43
$ws.Range("X43").Value = -14.36
$ws.Range("Y43").Value = 499.77
$ws.Range("W44").Value = 184.01
$ws.Range("X44").Value = 165.68
$ws.Range("Y44").Value = 709.05
$ws.Range("Z44").Value = 112.08
$ws.Range("W46").Value = -0.42
$ws.Range("X46").Value = 0.18
$ws.Range("Y46").Value = 0.63
$ws.Range("Z46").Value = 0.5600000000000001
$ws.Range("W48").Value = -24.43
$ws.Range("X48").Value = -14.36
$ws.Range("Y48").Value = 499.77
$ws.Range("W49").Value = 198.71
$ws.Range("X49").Value = 176.37
$ws.Range("Y49").Value = 719.86
$ws.Range("Z49").Value = 117.39
$ws.Range("W51").Value = 14.28
$ws.Range("X51").Value = 10.87
$ws.Range("Y51").Value = 11.44
$ws.Range("Z51").Value = 5.87
$ws.Range("W53").Value = -24.43
$ws.Range("X53").Value = -14.36
$ws.Range("Y53").Value = 499.77
$ws.Range("W54").Value = 185.06
$ws.Range("X54").Value = 171.06
$ws.Range("Y54").Value = 716.21
$ws.Range("Z54").Value = 115.93
$ws.Range("W56").Value = 0.63
$ws.Range("X56").Value = 5.56
$ws.Range("Y56").Value = 7.79
$ws.Range("Z56").Value = 4.41
$ws.Range("W58").Value = -24.43
$ws.Range("X58").Value = -14.36
$ws.Range("Y58").Value = 499.77
$ws.Range("W59").Value = 192.68
$ws.Range("X59").Value = 172.99
$ws.Range("Y59").Value = 717.5700000000001
$ws.Range("Z59").Value = 116.53
$ws.Range("W61").Value = 8.25
$ws.Range("X61").Value = 7.49
$ws.Range("Y61").Value = 9.15
$ws.Range("W63").Value = -24.43
$ws.Range("X63").Value = -14.36
$ws.Range("Y63").Value = 499.77
$ws.Range("W64").Value = 196.82
$ws.Range("X64").Value = 176.57
$ws.Range("Y64").Value = 721.27
$ws.Range("Z64").Value = 118.39
$ws.Range("W66").Value = 12.39
$ws.Range("X66").Value = 11.07
$ws.Range("Y66").Value = 12.85
$ws.Range("Z66").Value = 6.87
$ws.Range("W68").Value = -24.43
$ws.Range("X68").Value = -14.36
$ws.Range("Y68").Value = 499.77
$ws.Range("W69").Value = 198
$ws.Range("X69").Value = 178
$ws.Range("Y69").Value = 723.64
$ws.Range("Z69").Value = 119.79
$ws.Range("W71").Value = 13.57
$ws.Range("X71").Value = 12.5
$ws.Range("Y71").Value = 15.22
$ws.Range("Z71").Value = 8.27
$ws.Range("W73").Value = -24.43
$ws.Range("X73").Value = -14.36
$ws.Range("Y73").Value = 499.77
$ws.Range("W74").Value = 194.27
$ws.Range("X74").Value = 174.37
$ws.Range("Y74").Value = 718.71
$ws.Range("Z74").Value = 117.02
$ws.Range("W76").Value = 9.84
$ws.Range("X76").Value = 8.869999999999999
$ws.Range("Y76").Value = 10.29
$ws.Range("Z76").Value = 5.5
$ws.Range("W78").Value = -24.43
$ws.Range("X78").Value = -14.36
$ws.Range("Y78").Value = 499.77
$ws.Range("W79").Value = 195.55
$ws.Range("X79").Value = 175.42
$ws.Range("Y79").Value = 719.73
$ws.Range("Z79").Value = 117.65
$ws.Range("W81").Value = 11.12
$ws.Range("X81").Value = 9.92
$ws.Range("Y81").Value = 11.31
$ws.Range("Z81").Value = 6.12
$ws.Range("W83").Value = -24.43
$ws.Range("X83").Value = -14.36
$ws.Range("Y83").Value = 499.77
$ws.Range("W84").Value = 175.25
$ws.Range("X84").Value = 165.68
$ws.Range("Y84").Value = 710.53
$ws.Range("Z84").Value = 113.45
$ws.Range("W86").Value = -9.19
$ws.Range("X86").Value = 0.18
$ws.Range("Y86").Value = 2.11
$ws.Range("Z86").Value = 1.93
$ws.Range("W88").Value = -24.43
$ws.Range("X88").Value = -14.36
$ws.Range("Y88").Value = 499.77
$ws.Range("W89").Value = 110.95
$ws.Range("X89").Value = 94.62
$ws.Range("Y89").Value = 694.22
$ws.Range("Z89").Value = 106.01
$ws.Range("W90").Value = -54.32
$ws.Range("X90").Value = -56.17
$ws.Range("Y90").Value = 0
$ws.Range("W91").Value = -19.16
$ws.Range("X91").Value = -14.7
$ws.Range("Y91").Value = -14.2
$ws.Range("Z91").Value = -5.51
$ws.Range("W93").Value = -24.43
$ws.Range("X93").Value = -14.36
$ws.Range("Y93").Value = 499.77
